$wb = $excel.ActiveWorkbook

# The same table data is duplicated across the "展览" and "全部类型" sheets.
# Update the "想去人数" (interest count) column F for row 3 and row 5 on both.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 84
    $ws.Range("F5").Value = 25
}
